$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Yahoo")

# Copy header cell formatting (bold + border, existing style) into new H1 header cell
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Copy body cell formatting (border, existing style) into new H2:H4 data cells
$ws.Range("G2:G4").Copy()
$ws.Range("H2:H4").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Set the new column's values
$ws.Range("H1").Value = "Date"
$ws.Range("H2").Value = '${DATE}'
$ws.Range("H3").Value = '${DATE}'
$ws.Range("H4").Value = '${DATE}'

# Update the active selection to reflect where the user left off editing
$ws.Activate()
$ws.Range("H10").Select()
